# Apply the updated cryptocurrency price / 1h-volume snapshot to Sheet1.
# (cron-style data refresh: "Updated cryptos list ... with GitHub Actions")
#
# Price values in column D are stored as literal text in the source data
# (e.g. "42.442.02", "0.0960") to preserve exact formatting/trailing zeros.
# Values that look like plain numbers would otherwise be auto-converted to
# numeric by Excel on assignment, so those cells are briefly switched to
# Text format ("@") for the write and then restored to General.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.442.02"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "2.238.14"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.19"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.32"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.38"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +3.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0960"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.850"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "2.237.22"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "42.255.78"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("E18").Value = "  +9.37%  "
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.10"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.31"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +37.79%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  -4.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.77"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +3.71%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.29"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  +6.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.75"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.92"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.86"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +18.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0809"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -3.31%  "
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.93"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -8.77%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.48"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0309"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.27"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -7.23%  "
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("E40").Value = "  -4.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "63.56"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +4.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.202"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.92"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("E44").Value = "  -6.25%  "
$ws.Range("E45").Value = "  +2.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.996"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +3.77%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.14"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.08"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -2.37%  "
